$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new order row (row 6) with the waiter's completed order data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "dineIn"
$ws.Range("C6").Value = "[2, 2, 2, 2, 3]"
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = "InProgress"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
